# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The detail table (B16:G46) previously was grouped by period (column E)
# with each period listing all four/five workers. It is rebuilt here
# grouped by worker instead: each worker's own document number / name
# stays fixed down a contiguous block of rows, while the period (most
# recent first) varies. Along the way "DAIRO ENRIQUE VILLARREAL GOMEZ"'s
# Salario Basico (column G) drops from 2,000,000 to 1,000,000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tipo Doc, N Doc, Nombre, Periodo, Valor Mora, Salario Basico
$rows = @(
    @("CC", "45757080",   "JAQUELINE BUSTOS PASTRANA",      "2502", 50667, 2000000),
    @("CC", "45757080",   "JAQUELINE BUSTOS PASTRANA",      "2501", 80000, 2000000),
    @("CC", "45757080",   "JAQUELINE BUSTOS PASTRANA",      "2412", 80000, 2000000),
    @("CC", "45757080",   "JAQUELINE BUSTOS PASTRANA",      "2411", 80000, 2000000),
    @("CC", "45757080",   "JAQUELINE BUSTOS PASTRANA",      "2410", 80000, 2000000),
    @("CC", "45757080",   "JAQUELINE BUSTOS PASTRANA",      "2409", 80000, 2000000),
    @("CC", "45757080",   "JAQUELINE BUSTOS PASTRANA",      "2408", 80000, 2000000),
    @("CC", "78755471",   "JOSE GREGORIO ANAYA SANCHEZ",    "2410", 52000, 1300000),
    @("CC", "78755471",   "JOSE GREGORIO ANAYA SANCHEZ",    "2409", 52000, 1300000),
    @("CC", "78755471",   "JOSE GREGORIO ANAYA SANCHEZ",    "2408", 52000, 1300000),
    @("CC", "73181256",   "DAIRO ENRIQUE VILLARREAL GOMEZ", "2502", 50667, 1000000),
    @("CC", "73181256",   "DAIRO ENRIQUE VILLARREAL GOMEZ", "2501", 80000, 1000000),
    @("CC", "73181256",   "DAIRO ENRIQUE VILLARREAL GOMEZ", "2412", 80000, 1000000),
    @("CC", "73181256",   "DAIRO ENRIQUE VILLARREAL GOMEZ", "2411", 80000, 1000000),
    @("CC", "73181256",   "DAIRO ENRIQUE VILLARREAL GOMEZ", "2410", 80000, 1000000),
    @("CC", "73181256",   "DAIRO ENRIQUE VILLARREAL GOMEZ", "2409", 80000, 1000000),
    @("CC", "73181256",   "DAIRO ENRIQUE VILLARREAL GOMEZ", "2408", 80000, 1000000),
    @("CC", "1047390420", "DAIRO JOSE CARMONA NUNEZ",       "2502", 32933, 1300000),
    @("CC", "1047390420", "DAIRO JOSE CARMONA NUNEZ",       "2501", 52000, 1300000),
    @("CC", "1047390420", "DAIRO JOSE CARMONA NUNEZ",       "2412", 52000, 1300000),
    @("CC", "1047390420", "DAIRO JOSE CARMONA NUNEZ",       "2411", 52000, 1300000),
    @("CC", "1047390420", "DAIRO JOSE CARMONA NUNEZ",       "2410", 52000, 1300000),
    @("CC", "1047390420", "DAIRO JOSE CARMONA NUNEZ",       "2409", 52000, 1300000),
    @("CC", "1047390420", "DAIRO JOSE CARMONA NUNEZ",       "2408", 52000, 1300000),
    @("CC", "73187778",   "ALEXANDER ORTEGA MERCADO",       "2502", 32933, 1300000),
    @("CC", "73187778",   "ALEXANDER ORTEGA MERCADO",       "2501", 52000, 1300000),
    @("CC", "73187778",   "ALEXANDER ORTEGA MERCADO",       "2412", 52000, 1300000),
    @("CC", "73187778",   "ALEXANDER ORTEGA MERCADO",       "2411", 52000, 1300000),
    @("CC", "73187778",   "ALEXANDER ORTEGA MERCADO",       "2410", 52000, 1300000),
    @("CC", "73187778",   "ALEXANDER ORTEGA MERCADO",       "2409", 52000, 1300000),
    @("CC", "73187778",   "ALEXANDER ORTEGA MERCADO",       "2408", 52000, 1300000)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Range("B$r").Value = $data[0]
    $ws.Range("C$r").Value = $data[1]
    $ws.Range("D$r").Value = $data[2]
    $ws.Range("E$r").Value = $data[3]
    $ws.Range("F$r").Value = $data[4]
    $ws.Range("G$r").Value = $data[5]
}
